# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh values to the Gungnir Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 97.454544
$ws.Range("I9").Value = 81.71429000000001
$ws.Range("J9").Value = 125
$ws.Range("K9").Value = 81.71429000000001
$ws.Range("L9").Value = 125
$ws.Range("M9").Value = 87.28570999999999
$ws.Range("N9").Value = -463

$ws.Range("H74").Value = 4037
$ws.Range("I74").Value = 3906
$ws.Range("J74").Value = 4072.7273
$ws.Range("K74").Value = 3906
$ws.Range("L74").Value = 4072.7273
$ws.Range("M74").Value = -2970
$ws.Range("N74").Value = -5944.7273

$ws.Range("H77").Value = 4037
$ws.Range("I77").Value = 3906
$ws.Range("J77").Value = 4072.7273
$ws.Range("K77").Value = 19530
$ws.Range("L77").Value = 20363.6365
$ws.Range("M77").Value = -14850
$ws.Range("N77").Value = -29723.6365

$ws.Range("H86").Value = 11306.823
$ws.Range("I86").Value = 9215.846
$ws.Range("J86").Value = 12601.238
$ws.Range("K86").Value = 9215.846
$ws.Range("L86").Value = 12601.238
$ws.Range("M86").Value = -8092.846
$ws.Range("N86").Value = -14847.238

$ws.Range("H89").Value = 11306.823
$ws.Range("I89").Value = 9215.846
$ws.Range("J89").Value = 12601.238
$ws.Range("K89").Value = 46079.23
$ws.Range("L89").Value = 63006.19
$ws.Range("M89").Value = -40463.23
$ws.Range("N89").Value = -74238.19

$ws.Range("H137").Value = 1937.75
$ws.Range("I137").Value = 1499.8
$ws.Range("K137").Value = 4499.4
$ws.Range("M137").Value = -1949.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6063694
$ws.Range("I32").Value = 2781.6458
$ws.Range("J32").Value = 47624236
$ws.Range("K32").Value = 2781.6458
$ws.Range("L32").Value = 47624236
$ws.Range("M32").Value = -2494.6458
$ws.Range("N32").Value = -47624810

$ws.Range("H63").Value = 2492
$ws.Range("J63").Value = 2625.926
$ws.Range("L63").Value = 2625.926
$ws.Range("N63").Value = -3997.926

$ws.Range("H66").Value = 2492
$ws.Range("J66").Value = 2625.926
$ws.Range("L66").Value = 13129.63
$ws.Range("N66").Value = -19993.63

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 5381.8423
$ws.Range("J7").Value = 10085.4
$ws.Range("L7").Value = 10085.4
$ws.Range("N7").Value = -10311.4

$ws.Range("H31").Value = 4459.0347
$ws.Range("I31").Value = 900.5217
$ws.Range("J31").Value = 18100
$ws.Range("K31").Value = 900.5217
$ws.Range("L31").Value = 18100
$ws.Range("M31").Value = -605.5217
$ws.Range("N31").Value = -18690

$ws.Range("H34").Value = 4459.0347
$ws.Range("I34").Value = 900.5217
$ws.Range("J34").Value = 18100
$ws.Range("K34").Value = 900.5217
$ws.Range("L34").Value = 18100
$ws.Range("M34").Value = -698.5217
$ws.Range("N34").Value = -18504

$ws.Range("H132").Value = 14495353
$ws.Range("I132").Value = 2705.25
$ws.Range("J132").Value = 47621410
$ws.Range("K132").Value = 8115.75
$ws.Range("L132").Value = 142864230
$ws.Range("M132").Value = -5585.75
$ws.Range("N132").Value = -142869290

$ws.Range("H134").Value = 45456264
$ws.Range("I134").Value = 2050
$ws.Range("J134").Value = 71430100
$ws.Range("K134").Value = 6150
$ws.Range("L134").Value = 214290300
$ws.Range("M134").Value = -3615
$ws.Range("N134").Value = -214295370

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 94.77778000000001
$ws.Range("I6").Value = 71.85714
$ws.Range("J6").Value = 175
$ws.Range("K6").Value = 215.57142
$ws.Range("L6").Value = 525
$ws.Range("M6").Value = -102.57142
$ws.Range("N6").Value = -751

$ws.Range("H9").Value = 50001500
$ws.Range("J9").Value = 3000
$ws.Range("L9").Value = 9000
$ws.Range("N9").Value = -9448

$ws.Range("H17").Value = 206.2
$ws.Range("I17").Value = 206.2
$ws.Range("K17").Value = 618.5999999999999
$ws.Range("M17").Value = -449.5999999999999

$ws.Range("H114").Value = 670.5333000000001
$ws.Range("J114").Value = 1558
$ws.Range("L114").Value = 4674
$ws.Range("N114").Value = -11182

$ws.Range("H117").Value = 319.75
$ws.Range("I117").Value = 319.75
$ws.Range("K117").Value = 959.25
$ws.Range("M117").Value = 2482.75

$ws.Range("H121").Value = 14493344
$ws.Range("I121").Value = 264.25
$ws.Range("J121").Value = 22222986
$ws.Range("K121").Value = 792.75
$ws.Range("L121").Value = 66668958
$ws.Range("M121").Value = 517.25
$ws.Range("N121").Value = -66671578

$ws.Range("H131").Value = 860.6799999999999
$ws.Range("J131").Value = 908.76666
$ws.Range("L131").Value = 2726.29998
$ws.Range("N131").Value = -12806.29998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9045.904
$ws.Range("I70").Value = 10530.25
$ws.Range("K70").Value = 10530.25
$ws.Range("M70").Value = -10260.25

$ws.Range("H73").Value = 9045.904
$ws.Range("I73").Value = 10530.25
$ws.Range("K73").Value = 10530.25
$ws.Range("M73").Value = -9594.25

$ws.Range("H132").Value = 10442.929
$ws.Range("I132").Value = 3226.125
$ws.Range("J132").Value = 20065.334
$ws.Range("K132").Value = 9678.375
$ws.Range("L132").Value = 60196.00199999999
$ws.Range("M132").Value = -7148.375
$ws.Range("N132").Value = -65256.00199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4300
$ws.Range("I62").Value = 4300
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4300
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3676
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 4300
$ws.Range("I65").Value = 4300
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 21500
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -18380
$ws.Range("N65").ClearContents()

$ws.Range("H136").Value = 7248147.5
$ws.Range("I136").Value = 14707836
$ws.Range("J136").Value = 1592.2858
$ws.Range("K136").Value = 44123508
$ws.Range("L136").Value = 4776.857400000001
$ws.Range("M136").Value = -44120958
$ws.Range("N136").Value = -9876.857400000001
